# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 3, pushing the existing
# "Segunda" (row 3) and "Primera" (row 4) records down to rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (shifts old rows 3-4 -> 4-5)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's record
$ws.Cells.Item(3, 1).Value = 10
$ws.Cells.Item(3, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value = "La Araucanía"
$ws.Cells.Item(3, 4).Value = 44488
$ws.Cells.Item(3, 5).Value = 9
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100104
$ws.Cells.Item(3, 8).Value = "Frutos de pepita"
$ws.Cells.Item(3, 9).Value = 100104004
$ws.Cells.Item(3, 10).Value = "Níspero"
$ws.Cells.Item(3, 11).Value = "Californiana(o)"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 12000
$ws.Cells.Item(3, 15).Value = 12000
$ws.Cells.Item(3, 16).Value = 12000
$ws.Cells.Item(3, 17).Value = "$/bandeja 5 kilos"
$ws.Cells.Item(3, 18).Value = "La Ligua"
$ws.Cells.Item(3, 19).Value = 2400
$ws.Cells.Item(3, 20).Value = 5
